# Sprint 2 Retrospective and Sprint 3 Backlogs
# ------------------------------------------------------------
# This script reproduces the content edits made to "Sprint 2 Backlog.xlsx":
#   1. Row 1.4 description gains an explicit seed count ("Input 4 seeds
#      in each house" instead of "Input seeds in each house").
#   2. The placeholder "?" / "Nathaniel Leake" values that had been
#      filled into the "Team member Initial"/"Team member Actual"
#      columns for rows 1.3-1.14 (E7:F16) are cleared out again - those
#      tasks don't actually have an owner assigned yet.
#   3. Several rows that already had a "Team member Initial" (column E)
#      now also get a matching "Team member Actual" (column F) filled
#      in for rows 1.16, 1.17, and 2.2-2.9.
#   4. Row 3.8 (GUI / "display timer") gets its "Team member Actual"
#      (F44) filled in to match column E.
#   5. Leftover explicit-black-font styling (s="3") on the untouched
#      "Team member Actual" cells in the GUI section is cleared so those
#      cells fall back to the workbook's default cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update task description for row 8 (task 1.4)
$ws.Range("C8").Value = "Input 4 seeds in each house"

# 2. Clear out the placeholder "?" / "Nathaniel Leake" values in
#    E7:F16 (tasks 1.3 through 1.14 hadn't actually been assigned yet)
$ws.Range("E7:F16").Clear()

# 3. Fill in "Team member Actual" to match "Team member Initial"
#    for the remaining Sprint 1 rows and all of the Sprint 2 (AI) rows
$ws.Range("F20").Value = "Tony Huynh"
$ws.Range("F20").ClearFormats()
$ws.Range("F21").Value = "Tony Huynh"
$ws.Range("F21").ClearFormats()

$ws.Range("F25").Value = "Andrew Lam"
$ws.Range("F26").Value = "Tony Huynh"
$ws.Range("F27").Value = "Andrew Lam/Tony Huynh"
$ws.Range("F28").Value = "Tony Huynh"
$ws.Range("F29").Value = "Tony Huynh"
$ws.Range("F30").Value = "Andrew Lam"
$ws.Range("F31").Value = "Tony Huynh"
$ws.Range("F32").Value = "Tony Huynh"

# 4. Fill in "Team member Actual" for row 3.8 ("display timer")
$ws.Range("F44").Value = "Nathaniel Leake"

# 5. Drop the stale explicit-black-font style from the GUI section's
#    "Team member Actual" cells so they use the default cell style
$ws.Range("F37").ClearFormats()
$ws.Range("F38").ClearFormats()
$ws.Range("F40").ClearFormats()
$ws.Range("F42").ClearFormats()
$ws.Range("F43").ClearFormats()
